$d = $word.ActiveDocument

# Edit 1: merge "<id>" + "p020r_1" + "</id>" into a single run "<id>p020r_1</id>"
$d.Content.Find.Execute("<id>p020r_1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p020r_1</id>", 2)

# Edit 2: merge "<" + "i" + "d>" + "p020r_2" + "</id>" into a single run "<id>p020r_2</id>"
$d.Content.Find.Execute("<id>p020r_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p020r_2</id>", 2)
